$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script run found nothing new to report for row 120 (C120 used to
# hold the placeholder "NA" page number) - clear it out.
$ws.Range("C120").Value = ""

# Append the new result row produced by the latest script run.
# Force A121 to stay plain text (rather than being auto-parsed into a
# date serial) the same way the rest of the "Date" column already is,
# then drop the temporary text format so the cell keeps the sheet's
# default (unstyled) look, matching every other data row.
$ws.Range("A121").NumberFormat = "@"
$ws.Range("A121").Value = "2025-05-30"
$ws.Range("A121").Style = "Normal"

$ws.Range("B121").Value = "espèces exotiques envahissantes"
$ws.Range("C121").Value = 97
$ws.Range("D121").Value = 1
